$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: rename E1, add F1 ("Pressure") copying E1's bold/border style ---
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$ws.Range("E1").Value = "Feels Like"
$ws.Range("F1").Value = "Pressure"

# --- New "Pressure" values for F2:F73 ---
$fvals = @(1.626998490796339, 3.732457025590121, 3.248164495236074, 4.24611705247089, 3.288064461501717, 2.841543016677358, 0.4683738391404624, 1.838579015297168, 1.326011838408627, 1.219948216895418, 4.865052773762228, 1.965488623333802, 4.460232775885567, 3.155693129986314, 3.974056517708242, 2.51318546552596, 2.884519423131795, 2.46258846909432, 0.9762149389902225, 3.612260576307527, 1.403861812204279, 0.1215798321572692, 3.227361479535839, 0.8855533970352447, 4.702292921764571, 4.769642885012937, 4.574321951102243, 1.850793501277222, 0.07728308264433714, 4.641592812938627, 2.140920741586572, 4.833274095218348, 4.818099885446264, 4.265047277336801, 1.472244460347929, 1.925488643009626, 4.255683357584284, 1.584610025781388, 0.8474637334304624, 2.784006312291751, 4.680773870803905, 3.480148983374865, 2.850305850446825, 0.4858824688538427, 3.075036133495849, 4.950269250521316, 0.7004200761826201, 2.591648261818684, 4.386865359639777, 3.703843088771022, 3.48507870497634, 3.512420419935546, 1.797455756098776, 1.467959221322467, 4.046805777392568, 4.050566973395904, 4.335361592900519, 4.566202762782357, 2.556711994304689, 2.507581473435998, 3.991475894833876, 3.249819653888826, 3.509834386288517, 3.978963347180505, 4.450026709087831, 1.689975784257679, 1.87791476319972, 0.469909699204345, 2.89140070498087, 0.1797113689837104, 2.327990090662301, 2.713223173537883)

# --- Walk every data row: shift the Time (col B) forward by 2 days and write the Pressure value ---
for ($i = 0; $i -lt $fvals.Length; $i++) {
    $row = $i + 2
    $oldDate = $ws.Cells.Item($row, 2).Value2
    $ws.Cells.Item($row, 2).Value = $oldDate + 2
    $ws.Cells.Item($row, 6).Value = $fvals[$i]
}
